$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.363.62'
$ws.Range("E2").Value = '  +3.36%  '
$ws.Range("D3").Value = '3.616.04'
$ws.Range("E3").Value = '  +2.35%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '627.24'
$ws.Range("E5").Value = '  +3.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.43'
$ws.Range("E6").Value = '  +4.25%  '
$ws.Range("D7").Value = '3.615.73'
$ws.Range("E7").Value = '  +2.48%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.493'
$ws.Range("E9").Value = '  +1.96%  '
$ws.Range("E10").Value = '  +3.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.22'
$ws.Range("E11").Value = '  +5.77%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.440'
$ws.Range("E12").Value = '  +2.84%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000223'
$ws.Range("E13").Value = '  +0.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.18'
$ws.Range("E14").Value = '  +4.55%  '
$ws.Range("D15").Value = '4.230.47'
$ws.Range("E15").Value = '  +2.42%  '
$ws.Range("D16").Value = '3.621.00'
$ws.Range("E16").Value = '  +2.65%  '
$ws.Range("D17").Value = '69.587.85'
$ws.Range("E17").Value = '  +3.69%  '
$ws.Range("E18").Value = '  -0.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.61'
$ws.Range("E19").Value = '  +4.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.87'
$ws.Range("E20").Value = '  +3.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.12'
$ws.Range("E21").Value = '  +9.57%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '458.97'
$ws.Range("E22").Value = '  +3.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.640'
$ws.Range("E23").Value = '  +1.82%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.49'
$ws.Range("E24").Value = '  +0.90%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000136'
$ws.Range("E25").Value = '  +11.83%  '
$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.66'
$ws.Range("E26").Value = '  +4.46%  '
$ws.Range("B27").Value = 'WrappedeETH'
$ws.Range("C27").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D27").Value = '3.762.47'
$ws.Range("E27").Value = '  +2.34%  '
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.20'
$ws.Range("E29").Value = '  +11.50%  '
$ws.Range("E30").Value = '  +3.19%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.72'
$ws.Range("E31").Value = '  +2.72%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.176'
$ws.Range("E32").Value = '  +11.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.61'
$ws.Range("E33").Value = '  +7.62%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  -0.05%  '
$ws.Range("E35").Value = '  +5.54%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '26.38'
$ws.Range("E36").Value = '  +2.72%  '
$ws.Range("D37").Value = '3.603.28'
$ws.Range("E37").Value = '  +2.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.38'
$ws.Range("E38").Value = '  +4.80%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.37'
$ws.Range("E39").Value = '  +10.40%  '
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0926'
$ws.Range("E41").Value = '  +7.68%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.997'
$ws.Range("E42").Value = '  -0.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '176.17'
$ws.Range("E43").Value = '  +0.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.62'
$ws.Range("E44").Value = '  +1.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '31.78'
$ws.Range("E45").Value = '  +15.70%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.911'
$ws.Range("E46").Value = '  +2.36%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.38'
$ws.Range("E47").Value = '  +12.96%  '
$ws.Range("E48").Value = '  +8.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '46.30'
$ws.Range("E49").Value = '  +1.49%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.78'
$ws.Range("E50").Value = '  +2.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.267'
$ws.Range("E51").Value = '  +7.55%  '
